$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sending cluster sCs -> Target cluster ECs) ---
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.180598
$ws.Range("H2").Value = 18.541794
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 29.463391
$ws.Range("N2").Value = 88.390173
$ws.Range("O2").Value = 0.5524066141697463
$ws.Range("P2").Value = 0.5524066141697463
$ws.Range("Q2").Value = 182.101375487818
$ws.Range("R2").Value = 1638.912379390362
$ws.Range("S2").Value = 0.5524066141697463
$ws.Range("T2").Value = 0.5524066141697463

# --- Row 3 (Sending cluster sCs -> Target cluster FAPs) ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.180598
$ws.Range("H3").Value = 18.541794
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.136243
$ws.Range("N3").Value = 0.408729
$ws.Range("O3").Value = 0.002554408429577191
$ws.Range("P3").Value = 0.002554408429577191
$ws.Range("Q3").Value = 0.842063213314
$ws.Range("R3").Value = 7.578568919826
$ws.Range("S3").Value = 0.002554408429577191
$ws.Range("T3").Value = 0.002554408429577191

# --- Row 4 (Sending cluster sCs -> Target cluster M2) ---
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.180598
$ws.Range("H4").Value = 18.541794
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.897293666666666
$ws.Range("N4").Value = 11.691881
$ws.Range("O4").Value = 0.07307002777882997
$ws.Range("P4").Value = 0.07307002777882997
$ws.Range("Q4").Value = 24.08760544161266
$ws.Range("R4").Value = 216.788448974514
$ws.Range("S4").Value = 0.07307002777882997
$ws.Range("T4").Value = 0.07307002777882997

# --- Row 5 (new row: Sending cluster sCs -> Target cluster sCs) ---
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Edn3"
$ws.Range("C5").Value = "Ednrb"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.180598
$ws.Range("H5").Value = 18.541794
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 19.839492
$ws.Range("N5").Value = 59.518476
$ws.Range("O5").Value = 0.3719689496218466
$ws.Range("P5").Value = 0.3719689496218466
$ws.Range("Q5").Value = 122.619924576216
$ws.Range("R5").Value = 1103.579321185944
$ws.Range("S5").Value = 0.3719689496218466
$ws.Range("T5").Value = 0.3719689496218466
